$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '62.187.68'
Set-TextValue 'E2' '  -3.32%  '
Set-TextValue 'D3' '2.989.86'
Set-TextValue 'E3' '  -4.37%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '579.56'
Set-TextValue 'E5' '  -2.61%  '
Set-TextValue 'D6' '145.45'
Set-TextValue 'E6' '  -8.42%  '
Set-TextValue 'E7' '  -0.07%  '
Set-TextValue 'E8' '  -3.90%  '
Set-TextValue 'D9' '2.993.73'
Set-TextValue 'E9' '  -4.24%  '
Set-TextValue 'E10' '  -7.45%  '
Set-TextValue 'E11' '  -5.26%  '
Set-TextValue 'E12' '  -2.84%  '
Set-TextValue 'E13' '  -5.99%  '
Set-TextValue 'D14' '34.47'
Set-TextValue 'E14' '  -7.89%  '
Set-TextValue 'D15' '0.121'
Set-TextValue 'E15' '  +1.11%  '
Set-TextValue 'D16' '3.479.17'
Set-TextValue 'E16' '  -4.46%  '
Set-TextValue 'D17' '7.04'
Set-TextValue 'E17' '  -3.39%  '
Set-TextValue 'D18' '62.188.82'
Set-TextValue 'E18' '  -3.21%  '
Set-TextValue 'D19' '2.992.19'
Set-TextValue 'E19' '  -4.19%  '
Set-TextValue 'D20' '455.65'
Set-TextValue 'E20' '  -4.96%  '
Set-TextValue 'D21' '13.84'
Set-TextValue 'E21' '  -5.18%  '
Set-TextValue 'D22' '0.676'
Set-TextValue 'E22' '  -5.82%  '
Set-TextValue 'D23' '7.27'
Set-TextValue 'E23' '  -4.63%  '
Set-TextValue 'D24' '79.94'
Set-TextValue 'E24' '  -1.90%  '
Set-TextValue 'D25' '2.27'
Set-TextValue 'E25' '  -8.24%  '
Set-TextValue 'D26' '12.20'
Set-TextValue 'E26' '  -6.41%  '
Set-TextValue 'B27' 'Dai'
Set-TextValue 'C27' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D27' '0.998'
Set-TextValue 'E27' '  -0.21%  '
Set-TextValue 'B28' 'RenderToken'
Set-TextValue 'C28' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D28' '9.99'
Set-TextValue 'E28' '  -5.21%  '
Set-TextValue 'E29' '  +0.11%  '
Set-TextValue 'E30' '  -5.21%  '
Set-TextValue 'E31' '  -3.97%  '
Set-TextValue 'D32' '2.08'
Set-TextValue 'E32' '  -6.03%  '
Set-TextValue 'D33' '26.80'
Set-TextValue 'E33' '  -2.33%  '
Set-TextValue 'E34' '  -6.08%  '
Set-TextValue 'E35' '  -4.22%  '
Set-TextValue 'D36' '0.0₃0779'
Set-TextValue 'E36' '  -8.48%  '
Set-TextValue 'E37' '  -5.71%  '
Set-TextValue 'D38' '2.10'
Set-TextValue 'E38' '  -6.88%  '
Set-TextValue 'D39' '49.98'
Set-TextValue 'E39' '  -2.32%  '
Set-TextValue 'D40' '8.98'
Set-TextValue 'E40' '  -2.20%  '
Set-TextValue 'D41' '2.89'
Set-TextValue 'E41' '  -13.49%  '
Set-TextValue 'D42' '407.65'
Set-TextValue 'E42' '  -10.10%  '
Set-TextValue 'E43' '  -5.82%  '
Set-TextValue 'E44' '  -1.72%  '
Set-TextValue 'D45' '2.763.69'
Set-TextValue 'E45' '  -2.87%  '
Set-TextValue 'D46' '0.0349'
Set-TextValue 'E46' '  -4.86%  '
Set-TextValue 'D47' '38.44'
Set-TextValue 'E47' '  -4.32%  '
Set-TextValue 'D48' '127.35'
Set-TextValue 'E48' '  -2.67%  '
Set-TextValue 'E49' '  +0.05%  '
Set-TextValue 'E50' '  -2.91%  '
Set-TextValue 'D51' '23.64'
Set-TextValue 'E51' '  -8.40%  '
